$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Result column (C13:C21) flips from "f" to "t"
$ws.Range("C13:C21").Value = "t"

# Scroll the sheet view so row 5 is the top-left row, and move the active
# selection/cell to C11
$ws.Range("C11").Select()
$excel.ActiveWindow.ScrollRow = 5
